$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Append a new row (row 14) to the Card18 lookup table for the
# 151-300 tone range, crediting "محمد عبدالله" as the servicer.
# The leading apostrophes force text storage (matching the existing
# rows, which store every value - even numeric-looking ones - as
# text) without leaving the "number stored as text" quote-prefix
# marker behind once we reset the style back to Normal below.
$ws.Range("A14").Value = "'18"
$ws.Range("B14").Value = "'151"
$ws.Range("C14").Value = "'300"
$ws.Range("D14:N14").Value = "'"
$ws.Range("O14").Value = "محمد عبدالله"

# Clear the transient "quote prefix" number format picked up from the
# apostrophe-prefixed entries above so the new row's cell styling
# matches the rest of the sheet (no explicit style attribute).
$ws.Range("A14:N14").Style = "Normal"
